# laserresin4.xlsx - multiple changes to all parts
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row before row 84 first (pushes nema17_l/nema17_x from 84/85 to 85/86) ---
$ws.Rows.Item(84).Insert()

# --- row 83: drawer_z_from_bottom (=25) -> drawer_support_overlap (literal 40) ---
$ws.Range("A83").Value = "drawer_support_overlap"
$ws.Range("B83").Value = 40
$ws.Range("D83").Value = "new"

# --- row 84 (newly inserted, blank): drawer_support_from_top ---
$ws.Range("A84").Value = "drawer_support_from_top"
$ws.Range("B84").Formula = "=B38+B3+B1"
$ws.Range("D84").Value = "new"

# --- window_overlap (row31) / door_overlap (row35): literal 30 -> formula =B84 (now drawer_support_from_top) ---
$ws.Range("B31").Formula = "=B84"
$ws.Range("B35").Formula = "=B84"

# --- reservoir_z (row38): 50 -> 40 (cascades reservoir_x/y row36/37, bed_z2 row49) ---
$ws.Range("B38").Value = 40

# --- linear_bearing_spacing (row40): formula B36+66 -> B36+70 ---
$ws.Range("B40").Formula = "=B36+70"

# --- tray_border (row43): 20 -> 25, and now gets a D marker "y" (sm_thickness) ---
$ws.Range("B43").Value = 25
$ws.Range("D43").Value = "y"

# --- bed_x / bed_y (rows 46/47): formula B36-61/B37-61 -> B36-51/B37-51 ---
$ws.Range("B46").Formula = "=B36-51"
$ws.Range("B47").Formula = "=B37-51"

# --- view state: selection moved to B41 ---
$ws.Activate()
$ws.Range("B41").Select()
